# Update the cached "datetimeFigureOut" date field text wherever it
# appears (the slide master and every slide layout each carry their own
# copy of the placeholder) from 4/3/2023 -> 3/28/2024, and refresh the
# title of slide 1 from "UVOD" to "UVOD- 2024".

$p = $ppt.ActivePresentation

$oldDate = "4/3/2023"
$newDate = "3/28/2024"

# --- Slide Master -----------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every slide layout off the master ---------------------------------
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 1 title: "UVOD" -> "UVOD- 2024" ------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "UVOD") {
            $shape.TextFrame.TextRange.Text = "UVOD- 2024"
        }
    }
}
